{"js": "const replacements = [\n  [\"811\u00d72=\", \"940\u00d78=\"],\n  [\"749\u00d74=\", \"399\u00d72=\"],\n  [\"146\u00d79=\", \"958\u00d72=\"],\n  [\"832\u00d77=\", \"770\u00d73=\"],\n  [\"418\u00d75=\", \"512\u00d76=\"],\n  [\"253\u00d77=\", \"686\u00d74=\"],\n  [\"353\u00d75=\", \"454\u00d79=\"],\n  [\"788\u00d73=\", \"573\u00d79=\"],\n  [\"429\u00d77=\", \"649\u00d73=\"],\n  [\"616\u00d72=\", \"613\u00d79=\"],\n  [\"602\u00d79=\", \"814\u00d79=\"],\n  [\"168\u00d76=\", \"520\u00d74=\"],\n  [\"512\u00d79=\", \"926\u00d72=\"],\n  [\"664\u00d75=\", \"633\u00d73=\"],\n  [\"385\u00d74=\", \"713\u00d78=\"],\n  [\"673\u00d79=\", \"502\u00d79=\"],\n  [\"139\u00d79=\", \"120\u00d76=\"],\n  [\"444\u00d75=\", \"452\u00d75=\"],\n  [\"882\u00d75=\", \"133\u00d78=\"],\n  [\"498\u00d77=\", \"239\u00d73=\"],\n  [\"532\u00d78=\", \"648\u00d72=\"],\n  [\"682\u00d77=\", \"571\u00d77=\"],\n  [\"824\u00d77=\", \"477\u00d79=\"],\n  [\"389\u00d74=\", \"769\u00d75=\"],\n  [\"193\u00d78=\", \"145\u00d73=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"811\u00d72=\", \"940\u00d78=\"),\n    @(\"749\u00d74=\", \"399\u00d72=\"),\n    @(\"146\u00d79=\", \"958\u00d72=\"),\n    @(\"832\u00d77=\", \"770\u00d73=\"),\n    @(\"418\u00d75=\", \"512\u00d76=\"),\n    @(\"253\u00d77=\", \"686\u00d74=\"),\n    @(\"353\u00d75=\", \"454\u00d79=\"),\n    @(\"788\u00d73=\", \"573\u00d79=\"),\n    @(\"429\u00d77=\", \"649\u00d73=\"),\n    @(\"616\u00d72=\", \"613\u00d79=\"),\n    @(\"602\u00d79=\", \"814\u00d79=\"),\n    @(\"168\u00d76=\", \"520\u00d74=\"),\n    @(\"512\u00d79=\", \"926\u00d72=\"),\n    @(\"664\u00d75=\", \"633\u00d73=\"),\n    @(\"385\u00d74=\", \"713\u00d78=\"),\n    @(\"673\u00d79=\", \"502\u00d79=\"),\n    @(\"139\u00d79=\", \"120\u00d76=\"),\n    @(\"444\u00d75=\", \"452\u00d75=\"),\n    @(\"882\u00d75=\", \"133\u00d78=\"),\n    @(\"498\u00d77=\", \"239\u00d73=\"),\n    @(\"532\u00d78=\", \"648\u00d72=\"),\n    @(\"682\u00d77=\", \"571\u00d77=\"),\n    @(\"824\u00d77=\", \"477\u00d79=\"),\n    @(\"389\u00d74=\", \"769\u00d75=\"),\n    @(\"193\u00d78=\", \"145\u00d73=\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $pair[1], $wdReplaceAll)\n}"}
